$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepare a clean text cell holding the new date value "2024-03-18".
# Typing a date-shaped string into a cell normally gets auto-converted
# into a serial date by Excel, so we first force the source cell to
# Text format, assign the literal string, then reset its style back to
# Normal (so no stray formatting survives) before copying that text
# into every date cell in the table.
$dateCell = $ws.Cells.Item(2, 2)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2024-03-18"
$dateCell.Style = "Normal"

# Row 2 : No=1, Tanggal=2024-03-18, Jenis=cair, Instansi=Dekanat, Berat=1 Liter, Harga=10000
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 3).Value = "cair"
$ws.Cells.Item(2, 4).Value = "Dekanat"
$ws.Cells.Item(2, 5).Value = "1 Liter"
$ws.Cells.Item(2, 6).Value = 10000

# Row 3 : No=2, Tanggal=2024-03-18, Jenis=cair, Instansi=Dekanat, Berat=2 Liter, Harga=20000
$dateCell.Copy($ws.Cells.Item(3, 2)) | Out-Null
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 3).Value = "cair"
$ws.Cells.Item(3, 4).Value = "Dekanat"
$ws.Cells.Item(3, 5).Value = "2 Liter"
$ws.Cells.Item(3, 6).Value = 20000

# Row 4 : No=3, Tanggal=2024-03-18, Jenis=cair, Instansi=Dekanat, Berat=2 Liter, Harga=20000
$dateCell.Copy($ws.Cells.Item(4, 2)) | Out-Null
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 3).Value = "cair"
$ws.Cells.Item(4, 4).Value = "Dekanat"
$ws.Cells.Item(4, 5).Value = "2 Liter"
$ws.Cells.Item(4, 6).Value = 20000

# Row 5 : No=4, Tanggal=2024-03-18, Jenis=padat, Instansi=Dekanat, Berat=2 KG, Harga=30000
$dateCell.Copy($ws.Cells.Item(5, 2)) | Out-Null
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 3).Value = "padat"
$ws.Cells.Item(5, 4).Value = "Dekanat"
$ws.Cells.Item(5, 5).Value = "2 KG"
$ws.Cells.Item(5, 6).Value = 30000

# Row 6 : No=5, Tanggal=2024-03-18, Jenis=cair, Instansi=Dekanat, Berat=1 Liter, Harga=10000
$dateCell.Copy($ws.Cells.Item(6, 2)) | Out-Null
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 3).Value = "cair"
$ws.Cells.Item(6, 4).Value = "Dekanat"
$ws.Cells.Item(6, 5).Value = "1 Liter"
$ws.Cells.Item(6, 6).Value = 10000

# Row 7 : No=6, Tanggal=2024-03-18, Jenis=padat, Instansi=Dekanat, Berat=2 KG, Harga=30000
$dateCell.Copy($ws.Cells.Item(7, 2)) | Out-Null
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 3).Value = "padat"
$ws.Cells.Item(7, 4).Value = "Dekanat"
$ws.Cells.Item(7, 5).Value = "2 KG"
$ws.Cells.Item(7, 6).Value = 30000

# Row 8 (new) : No=7, Tanggal=2024-03-18, Jenis=cair, Instansi=Dekanat, Berat=1 Liter, Harga=10000
$dateCell.Copy($ws.Cells.Item(8, 2)) | Out-Null
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 3).Value = "cair"
$ws.Cells.Item(8, 4).Value = "Dekanat"
$ws.Cells.Item(8, 5).Value = "1 Liter"
$ws.Cells.Item(8, 6).Value = 10000

# Row 9 (new) : No=8, Tanggal=2024-03-18, Jenis=padat, Instansi=Dekanat, Berat=2 KG, Harga=30000
$dateCell.Copy($ws.Cells.Item(9, 2)) | Out-Null
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 3).Value = "padat"
$ws.Cells.Item(9, 4).Value = "Dekanat"
$ws.Cells.Item(9, 5).Value = "2 KG"
$ws.Cells.Item(9, 6).Value = 30000
